# Apply edit: "change order of surfactant pairs in DoE"
#
# 1. Rename the sheet.
# 2. Swap the surfactant-composition values (columns C:T) between the
#    row pairs (50,68) (51,69) ... (61,79) - i.e. samples 49-60 trade
#    their recorded surfactant-pair amounts with samples 67-78 while
#    keeping the ID/Sample-name columns (A,B) fixed to their row.
# 3. Fill in the previously-empty "Water" (U) and "Sample Density" (V)
#    measurements for rows 44-79.
# 4. Clear the stale direct cell formatting (a one-off font) that was
#    applied to a scattered set of C:T cells in rows 44-79.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "PhD_MasterDataset_OT_initial_Ja"

# --- Swap surfactant-pair values between row pairs ---------------------
$pairs = @(
    @(50,68), @(51,69), @(52,70), @(53,71), @(54,72), @(55,73),
    @(56,74), @(57,75), @(58,76), @(59,77), @(60,78), @(61,79)
)
foreach ($p in $pairs) {
    $r1 = $p[0]
    $r2 = $p[1]
    $addr1 = "C" + $r1 + ":T" + $r1
    $addr2 = "C" + $r2 + ":T" + $r2
    $v1 = $ws.Range($addr1).Value2
    $v2 = $ws.Range($addr2).Value2
    $ws.Range($addr1).Value2 = $v2
    $ws.Range($addr2).Value2 = $v1
}

# --- Fill in Water (U) / Sample Density (V) for rows 44-79 -------------
$uvdata = @(
    @(44, 75.7, 1.018),
    @(45, 77.5, 1.016),
    @(46, 74.5, 1.021),
    @(47, 68.6, 1.024),
    @(48, 75.6, 1.018),
    @(49, 74.1, 1.018),
    @(50, 72.8, 1.007),
    @(51, 71.6, 1.006),
    @(52, 73.8, 1.003),
    @(53, 73.6, 1.002),
    @(54, 78.6, 1.006),
    @(55, 70.6, 1.005),
    @(56, 74.6, 1.007),
    @(57, 75.9, 1.01),
    @(58, 73.6, 1.01),
    @(59, 71.8, 1.012),
    @(60, 75.4, 1.01),
    @(61, 74.8, 1.009),
    @(62, 78.9, 1.005),
    @(63, 75.1, 1.007),
    @(64, 77.4, 1.007),
    @(65, 72.5, 1.007),
    @(66, 76.7, 1.008),
    @(67, 71.5, 1.008),
    @(68, 71.3, 1.027),
    @(69, 70.9, 1.025),
    @(70, 75.8, 1.022),
    @(71, 76.2, 1.023),
    @(72, 71.5, 1.025),
    @(73, 72, 1.023),
    @(74, 73.6, 1.019),
    @(75, 72.2, 1.017),
    @(76, 71.3, 1.017),
    @(77, 77.9, 1.015),
    @(78, 72.3, 1.017),
    @(79, 74, 1.014)
)
foreach ($row in $uvdata) {
    $r = $row[0]
    $u = $row[1]
    $v = $row[2]
    $ws.Range("U" + $r).Value2 = $u
    $ws.Range("V" + $r).Value2 = $v
}

# --- Clear the leftover direct formatting on the C:T cells -------------
$ws.Range("C44:C49").ClearFormats()
$ws.Range("J44:J49").ClearFormats()
$ws.Range("D50:D55").ClearFormats()
$ws.Range("H50:H55").ClearFormats()
$ws.Range("E56:E61").ClearFormats()
$ws.Range("G56:G61").ClearFormats()
$ws.Range("L62:L67").ClearFormats()
$ws.Range("M62:M67").ClearFormats()
$ws.Range("F68:F73").ClearFormats()
$ws.Range("N68:N73").ClearFormats()
$ws.Range("I74:I79").ClearFormats()
$ws.Range("K74:K79").ClearFormats()
$ws.Range("O44:O79").ClearFormats()
$ws.Range("T44:T79").ClearFormats()

# --- Reset the view so the sheet opens scrolled to the top-left --------
$ws.Range("A1").Select()

Write-Output "edit complete"
